$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: read.php / Add functionality to show "In stock"/"Out of stock"
$ws.Range("A16").Value = "read.php"
$ws.Range("B16").Value = 'Add functionality to show "In stock"/"Out of stock"'

# Row 17: nav.php / Ensure it highlights the name of the current page in the navbar
$ws.Range("A17").Value = "nav.php"
$ws.Range("B17").Value = "Ensure it highlights the name of the current page in the navbar"
$ws.Range("D17").Value = $true
$ws.Range("E17").Value = 'Bootstrap has built-in functionality for this (the "active" class) so it took some PHP to insert a conditional statement that looks at the page name using the $_SERVER superglobal and inserts "active" into the class attribute if it matches.'

# Row 18: ALL / Get the background to work
$ws.Range("A18").Value = "ALL"
$ws.Range("B18").Value = "Get the background to work"

# Row heights for newly-filled wrapped-text rows
$ws.Rows.Item(16).RowHeight = 29
$ws.Rows.Item(17).RowHeight = 58

# Update the selection/view to match the post-edit state
$ws.Range("B18").Select() | Out-Null
